$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Diodes" (sheet14.xml): add a new row (row 3) for the RGB LED part.
# ---------------------------------------------------------------------------
$wsDiodes = $wb.Worksheets.Item("Diodes")

$wsDiodes.Range("A3").WrapText = $true
$wsDiodes.Range("A3").Value = "E6C0606RGBC3UDA"

$wsDiodes.Range("B3").WrapText = $true
$wsDiodes.Range("B3").Value = "LED E6C0606RGBC3UDA 0603 RGB 5V 20mA"

$wsDiodes.Range("C3").WrapText = $true
$wsDiodes.Range("C3").Value = "E6C0606RGBC3UDA"

$wsDiodes.Range("D3").WrapText = $true
$wsDiodes.Range("D3").Value = "'0603"

$wsDiodes.Range("E3").WrapText = $true
$wsDiodes.Range("E3").Value = "RGB"

$wsDiodes.Range("F3").WrapText = $true
$wsDiodes.Range("F3").Value = "20mA"

$wsDiodes.Range("G3").WrapText = $true
$wsDiodes.Range("G3").Value = "5V"

$wsDiodes.Range("H3").WrapText = $true
$wsDiodes.Range("H3").Value = "2.20V (max, red), 3.40V (max , green), 3.40V (max, blue)"

$wsDiodes.Range("I3").WrapText = $true
$wsDiodes.Range("I3").Value = "90 ℃"

$wsDiodes.Range("J3").WrapText = $true
$wsDiodes.Range("J3").Value = "EKINGLUX"

$wsDiodes.Range("K3").WrapText = $true
$wsDiodes.Range("K3").Value = "LED_0603"

$wsDiodes.Range("L3").WrapText = $true
$wsDiodes.Range("L3").Value = "Altium_Footprints.PcbLib"

$wsDiodes.Range("M3").WrapText = $true
$wsDiodes.Range("M3").Value = "LED_RGB_COMMON_P"

$wsDiodes.Range("N3").WrapText = $true
$wsDiodes.Range("N3").Value = "Altium_Schematic_Symbols.SchLib"

$wsDiodes.Range("O3").WrapText = $true
$wsDiodes.Range("O3").Value = "Datasheet"

$wsDiodes.Range("P3").WrapText = $true
$wsDiodes.Range("P3").Value = "https://datasheet.lcsc.com/lcsc/2203301730_EKINGLUX-E6C0606RGBC3UDA_C375569.pdf"

# Row heights: row 2 grows to fit the new (taller) content and row 3 matches it.
$wsDiodes.Rows.Item(2).RowHeight = 75
$wsDiodes.Rows.Item(3).RowHeight = 75

# Page setup (adds pageSetup/orientation to the sheet).
$wsDiodes.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Sheet "Connectors" (sheet5.xml): add two new rows (9 & 10) for the banana
# jack connectors (black & red).
# ---------------------------------------------------------------------------
$wsConn = $wb.Worksheets.Item("Connectors")

$wsConn.Range("A9").Value = "24.247.2"

$wsConn.Range("B9").WrapText = $true
$wsConn.Range("B9").Value = "JACK  24.247.2 TH 4mm  VERT BLACK"

$wsConn.Range("C9").Value = "24.247.2"

$wsConn.Range("D9").WrapText = $true
$wsConn.Range("D9").Value = "TH"

$wsConn.Range("E9").WrapText = $true
$wsConn.Range("E9").Value = "NA"

$wsConn.Range("F9").WrapText = $true
$wsConn.Range("F9").Value = "60V"

$wsConn.Range("G9").WrapText = $true
$wsConn.Range("G9").Value = "32A"

$wsConn.Range("H9").WrapText = $true
$wsConn.Range("H9").Value = "VERT BLACK BANANA JACK"

$wsConn.Range("I9").WrapText = $true
$wsConn.Range("I9").Value = "NA"

$wsConn.Range("J9").Value = "Changzhou Amass Elec"

$wsConn.Range("K9").WrapText = $true
$wsConn.Range("K9").Value = "BANANA_JACK_BLACK_24.247.2"

$wsConn.Range("L9").WrapText = $true
$wsConn.Range("L9").Value = "Altium_Footprints.PcbLib"

$wsConn.Range("M9").Value = "JACK_1P"

$wsConn.Range("N9").WrapText = $true
$wsConn.Range("N9").Value = "Altium_Schematic_Symbols.SchLib"

$wsConn.Range("O9").WrapText = $true
$wsConn.Range("O9").Value = "Datasheet"

$wsConn.Range("P9").WrapText = $true
$wsConn.Range("P9").Value = "https://www.lcsc.com/product-detail/Banana-Connectors-Alligator-Clips_Changzhou-Amass-Elec-24-247-2_C106272.html"

$wsConn.Range("A10").Value = "24.247.1"

$wsConn.Range("B10").WrapText = $true
$wsConn.Range("B10").Value = "JACK  24.247.1 TH 4mm  VERT RED"

$wsConn.Range("C10").Value = "24.247.1"

$wsConn.Range("D10").WrapText = $true
$wsConn.Range("D10").Value = "TH"

$wsConn.Range("E10").WrapText = $true
$wsConn.Range("E10").Value = "NA"

$wsConn.Range("F10").WrapText = $true
$wsConn.Range("F10").Value = "60V"

$wsConn.Range("G10").WrapText = $true
$wsConn.Range("G10").Value = "32A"

$wsConn.Range("H10").WrapText = $true
$wsConn.Range("H10").Value = "VERT RED BANANA JACK"

$wsConn.Range("I10").WrapText = $true
$wsConn.Range("I10").Value = "NA"

$wsConn.Range("J10").Value = "Changzhou Amass Elec"

$wsConn.Range("K10").WrapText = $true
$wsConn.Range("K10").Value = "BANANA_JACK_RED_24.247.1"

$wsConn.Range("L10").WrapText = $true
$wsConn.Range("L10").Value = "Altium_Footprints.PcbLib"

$wsConn.Range("M10").Value = "JACK_1P"

$wsConn.Range("N10").WrapText = $true
$wsConn.Range("N10").Value = "Altium_Schematic_Symbols.SchLib"

$wsConn.Range("O10").WrapText = $true
$wsConn.Range("O10").Value = "Datasheet"

$wsConn.Range("P10").WrapText = $true
$wsConn.Range("P10").Value = "https://www.lcsc.com/product-detail/Banana-Connectors-Alligator-Clips_Changzhou-Amass-Elec-24-247-2_C106272.html"

$wsConn.Rows.Item(9).RowHeight = 105
$wsConn.Rows.Item(10).RowHeight = 105

# ---------------------------------------------------------------------------
# Workbook view: make "Connectors" the active sheet/tab.
# ---------------------------------------------------------------------------
$wsConn.Activate()
$wsConn.Range("M3").Select()
